$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-like values (e.g. "27.128.81", "0.000008644") that
# must not be auto-converted to numbers, so force text format before assignment.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "27.128.81"
$ws.Cells.Item(2, 5).Value = "  -0.37%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.822.80"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.013"
$ws.Cells.Item(4, 5).Value = "  +0.18%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "312.16"
$ws.Cells.Item(5, 5).Value = "  -0.60%  "
$ws.Cells.Item(6, 5).Value = "  -0.18%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4616"
$ws.Cells.Item(7, 5).Value = "  -2.12%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3631"
$ws.Cells.Item(8, 5).Value = "  -1.71%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07295"
$ws.Cells.Item(9, 5).Value = "  -1.80%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.8707"
$ws.Cells.Item(10, 5).Value = "  -1.44%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "20.07"
$ws.Cells.Item(11, 5).Value = "  -1.93%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.878.35"
$ws.Cells.Item(12, 5).Value = "  +2.69%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.07631"
$ws.Cells.Item(13, 5).Value = "  +3.95%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.340"
$ws.Cells.Item(14, 5).Value = "  -2.56%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "92.37"
$ws.Cells.Item(15, 5).Value = "  -0.85%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.480"
$ws.Cells.Item(16, 5).Value = "  -1.47%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.011"
$ws.Cells.Item(17, 5).Value = "  -0.27%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000008644"
$ws.Cells.Item(18, 5).Value = "  -1.99%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.011"
$ws.Cells.Item(19, 5).Value = "  +0.00%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "27.416.40"
$ws.Cells.Item(20, 5).Value = "  +0.64%  "
$ws.Cells.Item(21, 5).Value = "  -2.28%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.209"
$ws.Cells.Item(22, 5).Value = "  -1.95%  "
$ws.Cells.Item(23, 5).Value = "  -1.37%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.096.68"
$ws.Cells.Item(24, 5).Value = "  +1.99%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.874"
$ws.Cells.Item(25, 5).Value = "  -1.33%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "151.38"
$ws.Cells.Item(26, 5).Value = "  -1.12%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.21"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.069"
$ws.Cells.Item(28, 5).Value = "  -4.78%  "
$ws.Cells.Item(29, 5).Value = "  -3.46%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "116.11"
$ws.Cells.Item(30, 5).Value = "  -1.57%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.08911"
$ws.Cells.Item(31, 5).Value = "  -0.18%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.961"
$ws.Cells.Item(32, 5).Value = "  +0.74%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.7340"
$ws.Cells.Item(33, 5).Value = "  -3.60%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.453"
$ws.Cells.Item(34, 5).Value = "  -2.20%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.135"
$ws.Cells.Item(35, 5).Value = "  -3.36%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.012"
$ws.Cells.Item(36, 5).Value = "  +0.15%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.460"
$ws.Cells.Item(37, 5).Value = "  +2.19%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.074"
$ws.Cells.Item(38, 5).Value = "  -2.76%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.05245"
$ws.Cells.Item(39, 5).Value = "  -1.88%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.01914"
$ws.Cells.Item(40, 5).Value = "  -2.42%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.930"
$ws.Cells.Item(41, 5).Value = "  -2.14%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "7.151"
$ws.Cells.Item(42, 5).Value = "  -2.48%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.5191"
$ws.Cells.Item(43, 5).Value = "  -3.01%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.1626"
$ws.Cells.Item(44, 5).Value = "  -2.43%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "8.275"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.4838"
$ws.Cells.Item(46, 5).Value = "  -2.46%  "
$ws.Cells.Item(47, 5).Value = "  -0.20%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "10.19"
$ws.Cells.Item(48, 5).Value = "  -3.12%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "103.34"
$ws.Cells.Item(49, 5).Value = "  -0.48%  "
$ws.Cells.Item(50, 5).Value = "  -2.41%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.06268"
$ws.Cells.Item(51, 5).Value = "  -0.82%  "
